$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (TC-001): trim the login step text and fix the redirect URL ---
$ws.Range("D2").Value = "输入用户名`"longchangkun`"`n输入密码`"Lck123456`"`n勾选我已阅读并同意复选框`n点击`"登录`"按钮`n"
$ws.Range("E2").Value = "跳转成功到页面，https://task-pre.renderbus.com/"

# --- Row 3 (TC-002): fill in the previously-empty test-steps cell ---
$ws.Range("D3").Value = "输入用户名`"longchangkun`"`n输入密码`"Lck123456`"`n勾选我已阅读并同意复选框`n点击`"登录`"按钮`n点击左侧导航栏的“统计”下面的 “云制作”"
$ws.Rows.Item(3).RowHeight = 84

# --- Update the active selection to D6 ---
$ws.Range("D6").Select()
